$d = $word.ActiveDocument

# --- 1) Rewrite the "Konsole" placeholder paragraph with real text ---
$d.Content.Find.Execute("(einfach Konsole?)", $true, $false, $false, $false, $false, $true, 1, $false, "Als Nutzerinterface benutzen wir die normale Konsolenausgabe. ", 2) | Out-Null

# --- 2) Fill in the following empty paragraph with the longer description ---
$p = $d.Paragraphs.Item(101)
$r = $p.Range
$body = $d.Range($r.Start, $r.End - 1)
$body.Text = "Zu Beginn des Spiels gibt es eine große Ausgabe des Spieltitels “FleetBattle“ mittels einer print Funktion, am Ende des Spiels gibt es eine große „Sieg“ Ausgabe in Verbindung mit dem jeweiligen Spielernamen, beziehungsweise eine „Niederlage“ Ausgabe, sollte gegen die CPU verloren worden sein. Unterstützt wird die Ausgabe noch durch farbliche Elemente. So wird zum Beispiel ein Schiff versenkt, wird die dazugehörige Ausgabe in grün ausgegeben."

# --- 3) Add the new ListLabel character styles (121-147) used by the class diagram's list ---
$newStyles = @(
  @{ num = 121; font = "Arial"; kind = "both" },
  @{ num = 122; font = "Courier New"; kind = "cs" },
  @{ num = 123; font = "Wingdings"; kind = "cs" },
  @{ num = 124; font = "Symbol"; kind = "cs" },
  @{ num = 125; font = "Courier New"; kind = "cs" },
  @{ num = 126; font = "Wingdings"; kind = "cs" },
  @{ num = 127; font = "Symbol"; kind = "cs" },
  @{ num = 128; font = "Courier New"; kind = "cs" },
  @{ num = 129; font = "Wingdings"; kind = "cs" },
  @{ num = 130; font = "Arial"; kind = "both" },
  @{ num = 131; font = "Courier New"; kind = "cs" },
  @{ num = 132; font = "Wingdings"; kind = "cs" },
  @{ num = 133; font = "Symbol"; kind = "cs" },
  @{ num = 134; font = "Courier New"; kind = "cs" },
  @{ num = 135; font = "Wingdings"; kind = "cs" },
  @{ num = 136; font = "Symbol"; kind = "cs" },
  @{ num = 137; font = "Courier New"; kind = "cs" },
  @{ num = 138; font = "Wingdings"; kind = "cs" },
  @{ num = 139; font = "Arial"; kind = "both" },
  @{ num = 140; font = "Courier New"; kind = "cs" },
  @{ num = 141; font = "Wingdings"; kind = "cs" },
  @{ num = 142; font = "Symbol"; kind = "cs" },
  @{ num = 143; font = "Courier New"; kind = "cs" },
  @{ num = 144; font = "Wingdings"; kind = "cs" },
  @{ num = 145; font = "Symbol"; kind = "cs" },
  @{ num = 146; font = "Courier New"; kind = "cs" },
  @{ num = 147; font = "Wingdings"; kind = "cs" }
)

foreach ($item in $newStyles) {
  $styleName = "ListLabel " + $item.num
  $s = $d.Styles.Add($styleName, 2)
  if ($item.kind -eq "both") {
    $s.Font.Name = $item.font
    $s.Font.NameBi = $item.font
  } else {
    $s.Font.NameBi = $item.font
  }
  $s.QuickStyle = $true
}

Write-Output "Done"
